$d = $word.ActiveDocument

# 1. Fix typo "iinterneta" -> "interneta" and drop the stray leading space
#    before "Ukoliko korisnik obavi kupovinu karte/karti putem interneta, ..."
$found1 = $d.Content.Find.Execute(
    " Ukoliko korisnik obavi kupovinu karte/karti putem iinterneta,",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ukoliko korisnik obavi kupovinu karte/karti putem interneta,", 2)
Write-Output "Fix typo/leading space: $found1"

# 2. Remove the trailing tab run after "... obanavljanem zaliha nekog suvenira."
$tabChar = [char]9
$searchText = "obanavljanem zaliha nekog suvenira." + $tabChar
$found2 = $d.Content.Find.Execute(
    $searchText, $false, $false, $false, $false, $false, $true, 1, $false,
    "obanavljanem zaliha nekog suvenira.", 2)
Write-Output "Remove trailing tab: $found2"
